
# Fruta / hortaliza, semanal
#
# A new weekly sample was inserted into the daily price log for
# "Terminal Hortofrutícola Agro Chillán - Pepino ensalada" at sheet
# row 102 (pushing every row from 102..148 down by one, to 103..149).
# The worksheet keeps its natural chronological/record order, so this
# is a true row insert (shifting the rest of the table), not just an
# overwrite of row 102's contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 102; everything below (old rows 102-148)
# shifts down to 103-149.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row with the new record. Columns A, B,
# C, E, F, G, H, I and R are constant across every data row in this
# sheet, so reuse the same values here.
$ws.Cells.Item(102, 1).Value  = 7
$ws.Cells.Item(102, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(102, 3).Value  = "Ñuble"
$ws.Cells.Item(102, 4).Value  = 44466
$ws.Cells.Item(102, 5).Value  = 16
$ws.Cells.Item(102, 6).Value  = 100112043
$ws.Cells.Item(102, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(102, 8).Value  = "Sin especificar"
$ws.Cells.Item(102, 9).Value  = "Primera"
$ws.Cells.Item(102, 10).Value = 300
$ws.Cells.Item(102, 11).Value = 16000
$ws.Cells.Item(102, 12).Value = 17000
$ws.Cells.Item(102, 13).Value = 16500
$ws.Cells.Item(102, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(102, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(102, 16).Value = 275
$ws.Cells.Item(102, 17).Value = 60
$ws.Cells.Item(102, 18).Value = "Hortaliza"
